$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912": header updates ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 06:58:31"
$ws1.Range("A3").Value = "Total filas: 60"

# ---- Sheet "LP1912-215": header update only ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:58:31"

# ---- Sheet "6203-6173": header updates ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:58:31"
$ws3.Range("A3").Value = "Total filas: 14"

# ---- Sheet "LP1912": rows 35-65 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A35").Value = "06:58:31"
$ws1.Range("B35").Value = "07:00"
$ws1.Range("C35").Value = "14_ABASTO"
$ws1.Range("D35").Value = 2
$ws1.Range("E35").Value = "LP1912"

$ws1.Range("A36").Value = "06:15:33"
$ws1.Range("B36").Value = "07:01"
$ws1.Range("C36").Value = "16_SANTA ANA"
$ws1.Range("D36").Value = 46
$ws1.Range("E36").Value = "LP1912"

$ws1.Range("A37").Value = "05:53:46"
$ws1.Range("B37").Value = "07:04"
$ws1.Range("C37").Value = "23_HERNANDEZ"
$ws1.Range("D37").Value = 71
$ws1.Range("E37").Value = "LP1912"

$ws1.Range("A38").Value = "05:18:42"
$ws1.Range("B38").Value = "07:05"
$ws1.Range("C38").Value = "15_ABASTO"
$ws1.Range("D38").Value = 107
$ws1.Range("E38").Value = "LP1912"

$ws1.Range("A39").Value = "06:44:40"
$ws1.Range("B39").Value = "07:05"
$ws1.Range("C39").Value = "23_HERNANDEZ"
$ws1.Range("D39").Value = 21
$ws1.Range("E39").Value = "LP1912"

$ws1.Range("A40").Value = "05:18:42"
$ws1.Range("B40").Value = "07:07"
$ws1.Range("C40").Value = "225_GOMEZ"
$ws1.Range("D40").Value = 109
$ws1.Range("E40").Value = "LP1912"

$ws1.Range("A41").Value = "06:58:31"
$ws1.Range("B41").Value = "07:08"
$ws1.Range("C41").Value = "15_ABASTO"
$ws1.Range("D41").Value = 10
$ws1.Range("E41").Value = "LP1912"

$ws1.Range("A42").Value = "06:44:40"
$ws1.Range("B42").Value = "07:09"
$ws1.Range("C42").Value = "15_ABASTO"
$ws1.Range("D42").Value = 25
$ws1.Range("E42").Value = "LP1912"

$ws1.Range("A43").Value = "05:18:42"
$ws1.Range("B43").Value = "07:11"
$ws1.Range("C43").Value = "215A_EL PATO"
$ws1.Range("D43").Value = 113
$ws1.Range("E43").Value = "LP1912"

$ws1.Range("A44").Value = "05:18:42"
$ws1.Range("B44").Value = "07:15"
$ws1.Range("C44").Value = "11_ETCHEVERRY"
$ws1.Range("D44").Value = 117
$ws1.Range("E44").Value = "LP1912"

$ws1.Range("A45").Value = "06:44:40"
$ws1.Range("B45").Value = "07:16"
$ws1.Range("C45").Value = "16_SANTA ANA"
$ws1.Range("D45").Value = 32
$ws1.Range("E45").Value = "LP1912"

$ws1.Range("A46").Value = "05:53:46"
$ws1.Range("B46").Value = "07:21"
$ws1.Range("C46").Value = "26_HERNANDEZ"
$ws1.Range("D46").Value = 88
$ws1.Range("E46").Value = "LP1912"

$ws1.Range("A47").Value = "06:15:33"
$ws1.Range("B47").Value = "07:23"
$ws1.Range("C47").Value = "10_OLMOS"
$ws1.Range("D47").Value = 68
$ws1.Range("E47").Value = "LP1912"

$ws1.Range("A48").Value = "05:53:46"
$ws1.Range("B48").Value = "07:31"
$ws1.Range("C48").Value = "11_ETCHEVERRY"
$ws1.Range("D48").Value = 98
$ws1.Range("E48").Value = "LP1912"

$ws1.Range("A49").Value = "05:53:46"
$ws1.Range("B49").Value = "07:32"
$ws1.Range("C49").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D49").Value = 99
$ws1.Range("E49").Value = "LP1912"

$ws1.Range("A50").Value = "05:53:46"
$ws1.Range("B50").Value = "07:36"
$ws1.Range("C50").Value = "27_EL RETIRO"
$ws1.Range("D50").Value = 103
$ws1.Range("E50").Value = "LP1912"

$ws1.Range("A51").Value = "05:53:46"
$ws1.Range("B51").Value = "07:39"
$ws1.Range("C51").Value = "10_OLMOS"
$ws1.Range("D51").Value = 106
$ws1.Range("E51").Value = "LP1912"

$ws1.Range("A52").Value = "05:53:46"
$ws1.Range("B52").Value = "07:47"
$ws1.Range("C52").Value = "14_ABASTO"
$ws1.Range("D52").Value = 114
$ws1.Range("E52").Value = "LP1912"

$ws1.Range("A53").Value = "05:53:46"
$ws1.Range("B53").Value = "07:51"
$ws1.Range("C53").Value = "215D_EL PATO"
$ws1.Range("D53").Value = 118
$ws1.Range("E53").Value = "LP1912"

$ws1.Range("A54").Value = "06:58:31"
$ws1.Range("B54").Value = "07:58"
$ws1.Range("C54").Value = "16_SANTA ANA"
$ws1.Range("D54").Value = 60
$ws1.Range("E54").Value = "LP1912"

$ws1.Range("A55").Value = "06:15:33"
$ws1.Range("B55").Value = "08:07"
$ws1.Range("C55").Value = "16_SANTA ANA"
$ws1.Range("D55").Value = 112
$ws1.Range("E55").Value = "LP1912"

$ws1.Range("A56").Value = "06:15:33"
$ws1.Range("B56").Value = "08:12"
$ws1.Range("C56").Value = "15_ABASTO"
$ws1.Range("D56").Value = 117
$ws1.Range("E56").Value = "LP1912"

$ws1.Range("A57").Value = "06:58:31"
$ws1.Range("B57").Value = "08:13"
$ws1.Range("C57").Value = "10_OLMOS"
$ws1.Range("D57").Value = 75
$ws1.Range("E57").Value = "LP1912"

$ws1.Range("A58").Value = "06:44:40"
$ws1.Range("B58").Value = "08:21"
$ws1.Range("C58").Value = "26_HERNANDEZ"
$ws1.Range("D58").Value = 97
$ws1.Range("E58").Value = "LP1912"

$ws1.Range("A59").Value = "06:44:40"
$ws1.Range("B59").Value = "08:22"
$ws1.Range("C59").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D59").Value = 98
$ws1.Range("E59").Value = "LP1912"

$ws1.Range("A60").Value = "06:44:40"
$ws1.Range("B60").Value = "08:23"
$ws1.Range("C60").Value = "215B_EL PATO"
$ws1.Range("D60").Value = 99
$ws1.Range("E60").Value = "LP1912"

$ws1.Range("A61").Value = "06:44:40"
$ws1.Range("B61").Value = "08:27"
$ws1.Range("C61").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D61").Value = 103
$ws1.Range("E61").Value = "LP1912"

$ws1.Range("A62").Value = "06:58:31"
$ws1.Range("B62").Value = "08:34"
$ws1.Range("C62").Value = "23_HERNANDEZ"
$ws1.Range("D62").Value = 96
$ws1.Range("E62").Value = "LP1912"

$ws1.Range("A63").Value = "06:44:40"
$ws1.Range("B63").Value = "08:35"
$ws1.Range("C63").Value = "23_HERNANDEZ"
$ws1.Range("D63").Value = 111
$ws1.Range("E63").Value = "LP1912"

$ws1.Range("A64").Value = "06:44:40"
$ws1.Range("B64").Value = "08:42"
$ws1.Range("C64").Value = "81_EL PELIGRO"
$ws1.Range("D64").Value = 118
$ws1.Range("E64").Value = "LP1912"

$ws1.Range("A65").Value = "06:58:31"
$ws1.Range("B65").Value = "08:54"
$ws1.Range("C65").Value = "17_ROMERO"
$ws1.Range("D65").Value = 116
$ws1.Range("E65").Value = "LP1912"

# ---- Sheet "6203-6173": rows 17-19 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A17").Value = "06:58:31"
$ws3.Range("B17").Value = "08:07"
$ws3.Range("C17").Value = "215C_LA PLATA"
$ws3.Range("D17").Value = 69
$ws3.Range("E17").Value = "L6203"

$ws3.Range("A18").Value = "06:58:31"
$ws3.Range("B18").Value = "08:37"
$ws3.Range("C18").Value = "215A_LA PLATA"
$ws3.Range("D18").Value = 99
$ws3.Range("E18").Value = "L6173"

$ws3.Range("A19").Value = "06:44:40"
$ws3.Range("B19").Value = "08:41"
$ws3.Range("C19").Value = "215A_LA PLATA"
$ws3.Range("D19").Value = 117
$ws3.Range("E19").Value = "L6173"

